$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "11.3"
$ws.Range("D2").Value = "μmol/L"
$ws.Range("E2").Value = "2-25"

$ws.Range("D3").Value = "μmol/L"
$ws.Range("E3").Value = "0-5"

$ws.Range("D4").Value = "μmol/L"
$ws.Range("E4").Value = "0-20"

$ws.Range("E5").Value = "62-85"

$ws.Range("E6").Value = "35-55"

$ws.Range("E7").Value = "20-40"

$ws.Range("E8").Value = "1.2-2.4"

$ws.Range("B9").Value = "白蛋白"
$ws.Range("E9").Value = "200-400"

$ws.Range("E10").Value = "7-40"

$ws.Range("E11").Value = "13-35"

$ws.Range("E13").Value = "35-135"

$ws.Range("B14").Value = "γ谷氨酰转肽酶"
$ws.Range("E14").Value = "7-45"

$ws.Range("D15").Value = "μmol/L"
$ws.Range("E15").Value = "0-12"

$ws.Range("E16").Value = "4000-12000"

$ws.Range("D17").Value = "mmol/L"
$ws.Range("E17").Value = "2.5-6.4"

$ws.Range("D18").Value = "μmol/L"
$ws.Range("E18").Value = "40-97"

$ws.Range("D19").Value = "μmol/L"
$ws.Range("E19").Value = "150-430"

$ws.Range("A20").Value = "Cys"
$ws.Range("D20").Value = "g/L"
$ws.Range("E20").Value = "0-1.16"

$ws.Range("D21").Value = "mg/L"
$ws.Range("E21").Value = "25-70"

$ws.Range("A22").Value = "C"
$ws.Range("B22").Value = "总二氧化碳"
$ws.Range("D22").Value = "mmol/L"
$ws.Range("E22").Value = "20-30"
